$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four distinct rows of the repeating pattern already present
# in rows 2-5 / 6-9 / 10-13 are appended again as rows 14-17.
$data = @(
    @(1.1993322807698887, 2.3758877717639884, 2.413486364972186, -0.86585001746396684, -0.78012023365383742, 0, 8.3366194783562833, 0.86585001746396684, 4.5494284741316866),
    @([double]"2.688821387764051e-17", 2.7925609058034806, $null, -0.010657853425638181, 0.031468739706286171, [double]"3.7665825361947448e+18", [double]"1.275668646441314e+18", 0.031468739706286171, 0.23430532913982935),
    @(0, 2.7122120396162424, $null, -0.97942001697994874, 0.61367412211482841, 65535, 65535, 0.97942001697994874, 0.53369547192961131),
    @(0, 2.1421186466279774, $null, -0.099589484540936551, 0.22536832624391215, 65535, 65535, 0.22536832624391215, 0.07159329501494785)
)

$startRow = 14
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $val = $values[$j]
        $cell = $ws.Cells.Item($row, $col)
        if ($null -ne $val) {
            $cell.Value = $val
        } else {
            # Column C is left blank on these rows (matches rows 3/4/5 pattern),
            # but the cell itself still needs to materialize in the sheet.
            $cell.Font.Size = 11
        }
    }
}
